# Refresh the "cryptos" price/volume table with the latest scraped values
# (GitHub Actions cron update). Updates Price (D) and Volume(1h) (E) for
# every coin row, plus a rank swap between EnergySwap and BabyDogeCoin
# (rows 48/49) whose market caps crossed over in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain text formatting so numeric-looking strings
# (e.g. "0.9990", "29.339.81") are not auto-coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.339.81'
$ws.Range("E2").Value = '  -0.29%  '

$ws.Range("D3").Value = '1.842.83'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '240.52'
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").Value = '0.6285'

$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").Value = '0.07465'
$ws.Range("E8").Value = '  -2.38%  '

$ws.Range("D9").Value = '0.2893'
$ws.Range("E9").Value = '  -0.65%  '

$ws.Range("D10").Value = '24.34'
$ws.Range("E10").Value = '  -2.36%  '

$ws.Range("D11").Value = '0.07727'
$ws.Range("E11").Value = '  -0.20%  '

$ws.Range("D12").Value = '1.843.02'
$ws.Range("E12").Value = '  -2.39%  '

$ws.Range("E13").Value = '  -0.76%  '

$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("D15").Value = '0.00001016'
$ws.Range("E15").Value = '  -4.42%  '

$ws.Range("D16").Value = '82.03'
$ws.Range("E16").Value = '  -1.67%  '

$ws.Range("D17").Value = '6.119'

$ws.Range("D18").Value = '29.369.77'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").Value = '228.14'
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").Value = '0.9996'
$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").Value = '7.435'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").Value = '0.9995'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").Value = '158.65'
$ws.Range("E24").Value = '  +0.67%  '

$ws.Range("E25").Value = '  -0.88%  '

$ws.Range("D26").Value = '8.410'
$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("D27").Value = '17.55'

$ws.Range("D28").Value = '0.06484'
$ws.Range("E28").Value = '  +15.49%  '

$ws.Range("E29").Value = '  +0.37%  '

$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("E31").Value = '  -1.19%  '

$ws.Range("D32").Value = '4.044'
$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("D33").Value = '1.820'
$ws.Range("E33").Value = '  -1.47%  '

$ws.Range("E34").Value = '  -2.08%  '

$ws.Range("D35").Value = '0.6958'
$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("D36").Value = '2.583'
$ws.Range("E36").Value = '  -0.34%  '

$ws.Range("D37").Value = '1.260.51'
$ws.Range("E37").Value = '  +2.52%  '

$ws.Range("E38").Value = '  +4.18%  '

$ws.Range("D39").Value = '0.01811'
$ws.Range("E39").Value = '  +0.38%  '

$ws.Range("D40").Value = '6.522'
$ws.Range("E40").Value = '  +1.20%  '

$ws.Range("D41").Value = '0.9111'
$ws.Range("E41").Value = '  +0.30%  '

$ws.Range("D42").Value = '0.9987'
$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("D43").Value = '2.003.96'
$ws.Range("E43").Value = '  -12.82%  '

$ws.Range("D44").Value = '101.24'
$ws.Range("E44").Value = '  -0.64%  '

$ws.Range("D45").Value = '66.25'
$ws.Range("E45").Value = '  +0.22%  '

$ws.Range("E46").Value = '  +1.25%  '

$ws.Range("D47").Value = '7.014'
$ws.Range("E47").Value = '  -2.50%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.00000000115'
$ws.Range("E48").Value = '  -3.91%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.042'
$ws.Range("E49").Value = '  +0.55%  '

$ws.Range("E50").Value = '  -2.06%  '

$ws.Range("D51").Value = '1.674'
$ws.Range("E51").Value = '  -0.71%  '
